$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# regen save_data to use K instead of Strike#: overwrite column G (K) values
$newK = @{
    2  = 4
    3  = 4
    4  = 4
    5  = 3
    6  = 6
    7  = 4
    8  = 5
    9  = 5
    10 = 4
    11 = 5
    12 = 3
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
